$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 156, shifting existing rows 156-183 down to 157-184
$ws.Rows.Item(156).Insert()

# Populate the newly inserted row 156 with the new record's data
$ws.Cells.Item(156, 1).Value = 9
$ws.Cells.Item(156, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(156, 3).Value = "Metropolitana"
$ws.Cells.Item(156, 4).Value = 44522
$ws.Cells.Item(156, 5).Value = 13
$ws.Cells.Item(156, 6).Value = "Fruta"
$ws.Cells.Item(156, 7).Value = 100103
$ws.Cells.Item(156, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(156, 9).Value = 100103001
$ws.Cells.Item(156, 10).Value = "Cereza"
$ws.Cells.Item(156, 11).Value = "Early Burlat"
$ws.Cells.Item(156, 12).Value = "Primera"
$ws.Cells.Item(156, 13).Value = 350
$ws.Cells.Item(156, 14).Value = 28000
$ws.Cells.Item(156, 15).Value = 28000
$ws.Cells.Item(156, 16).Value = 28000
$ws.Cells.Item(156, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(156, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(156, 19).Value = 2800
$ws.Cells.Item(156, 20).Value = 10
